$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Itgb8 -> ECs (D2 changes from FAPs(21) to ECs(20)), plus new G..T values
$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 12.056684
$ws.Range("H2").Value = 36.170052
$ws.Range("I2").Value = 0.06307822458376462
$ws.Range("J2").Value = 0.06307822458376462
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.027767
$ws.Range("N2").Value = 0.083301
$ws.Range("O2").Value = 0.002923627791763407
$ws.Range("P2").Value = 0.002923627791763407
$ws.Range("Q2").Value = 0.334777944628
$ws.Range("R2").Value = 3.013001501652
$ws.Range("S2").Value = 0.000184417250448188
$ws.Range("T2").Value = 0.000184417250448188

# Row 3: ECs -> Itgb8 -> FAPs (D3 changes from MuSCs(22) to FAPs(21))
$ws.Range("D3").Value = "FAPs"
$ws.Range("G3").Value = 12.056684
$ws.Range("H3").Value = 36.170052
$ws.Range("I3").Value = 0.06307822458376462
$ws.Range("J3").Value = 0.06307822458376462
$ws.Range("M3").Value = 4.237840333333334
$ws.Range("N3").Value = 12.713521
$ws.Range("O3").Value = 0.4462083687682946
$ws.Range("P3").Value = 0.4462083687682946
$ws.Range("Q3").Value = 51.09430174145466
$ws.Range("R3").Value = 459.848715673092
$ws.Range("S3").Value = 0.02814603169632175
$ws.Range("T3").Value = 0.02814603169632175

# Row 4: FAPs -> Itgb8 -> FAPs becomes ECs -> Itgb8 -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 12.056684
$ws.Range("H4").Value = 36.170052
$ws.Range("I4").Value = 0.06307822458376462
$ws.Range("J4").Value = 0.06307822458376462
$ws.Range("M4").Value = 5.231839666666667
$ws.Range("N4").Value = 15.695519
$ws.Range("O4").Value = 0.5508680034399419
$ws.Range("P4").Value = 0.5508680034399419
$ws.Range("Q4").Value = 63.07863759966533
$ws.Range("R4").Value = 567.707738396988
$ws.Range("S4").Value = 0.03474777563699468
$ws.Range("T4").Value = 0.03474777563699468

# Row 5: FAPs -> Itgb8 -> MuSCs becomes FAPs -> Itgb8 -> ECs
$ws.Range("D5").Value = "ECs"
$ws.Range("I5").Value = 0.1315309049843414
$ws.Range("J5").Value = 0.1315309049843414
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.027767
$ws.Range("N5").Value = 0.083301
$ws.Range("O5").Value = 0.002923627791763407
$ws.Range("P5").Value = 0.002923627791763407
$ws.Range("Q5").Value = 0.6980799842780001
$ws.Range("R5").Value = 6.282719858502
$ws.Range("S5").Value = 0.0003845474092880125
$ws.Range("T5").Value = 0.0003845474092880125

# Row 6: MuSCs -> Itgb8 -> FAPs becomes FAPs -> Itgb8 -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("G6").Value = 25.140634
$ws.Range("H6").Value = 75.421902
$ws.Range("I6").Value = 0.1315309049843414
$ws.Range("J6").Value = 0.1315309049843414
$ws.Range("O6").Value = 0.4462083687682946
$ws.Range("P6").Value = 0.4462083687682946
$ws.Range("Q6").Value = 106.5419927707713
$ws.Range("R6").Value = 958.8779349369421
$ws.Range("S6").Value = 0.05869019055568051
$ws.Range("T6").Value = 0.0586901905556805

# Row 7: MuSCs -> Itgb8 -> MuSCs becomes FAPs -> Itgb8 -> MuSCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("G7").Value = 25.140634
$ws.Range("H7").Value = 75.421902
$ws.Range("I7").Value = 0.1315309049843414
$ws.Range("J7").Value = 0.1315309049843414
$ws.Range("M7").Value = 5.231839666666667
$ws.Range("O7").Value = 0.5508680034399419
$ws.Range("P7").Value = 0.5508680034399419
$ws.Range("Q7").Value = 131.5317662063487
$ws.Range("R7").Value = 1183.785895857138
$ws.Range("S7").Value = 0.07245616701937282
$ws.Range("T7").Value = 0.07245616701937282

# New rows 8,9,10: MuSCs -> Vtn -> Itgb8 -> ECs/FAPs/MuSCs
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Vtn"
$ws.Range("C8").Value = "Itgb8"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 153.9412893333333
$ws.Range("H8").Value = 461.8238680000001
$ws.Range("I8").Value = 0.8053908704318941
$ws.Range("J8").Value = 0.8053908704318941
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.027767
$ws.Range("N8").Value = 0.083301
$ws.Range("O8").Value = 0.002923627791763407
$ws.Range("P8").Value = 0.002923627791763407
$ws.Range("Q8").Value = 4.274487780918667
$ws.Range("R8").Value = 38.47039002826801
$ws.Range("S8").Value = 0.002354663132027207
$ws.Range("T8").Value = 0.002354663132027207

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Vtn"
$ws.Range("C9").Value = "Itgb8"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 153.9412893333333
$ws.Range("H9").Value = 461.8238680000001
$ws.Range("I9").Value = 0.8053908704318941
$ws.Range("J9").Value = 0.8053908704318941
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.237840333333334
$ws.Range("N9").Value = 12.713521
$ws.Range("O9").Value = 0.4462083687682946
$ws.Range("P9").Value = 0.4462083687682946
$ws.Range("Q9").Value = 652.3786049021365
$ws.Range("R9").Value = 5871.407444119229
$ws.Range("S9").Value = 0.3593721465162923
$ws.Range("T9").Value = 0.3593721465162923

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Vtn"
$ws.Range("C10").Value = "Itgb8"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 153.9412893333333
$ws.Range("H10").Value = 461.8238680000001
$ws.Range("I10").Value = 0.8053908704318941
$ws.Range("J10").Value = 0.8053908704318941
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.231839666666667
$ws.Range("N10").Value = 15.695519
$ws.Range("O10").Value = 0.5508680034399419
$ws.Range("P10").Value = 0.5508680034399419
$ws.Range("Q10").Value = 805.3961438719437
$ws.Range("R10").Value = 7248.565294847494
$ws.Range("S10").Value = 0.4436640607835744
$ws.Range("T10").Value = 0.4436640607835744
